$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new "Week 45" column header in AT1
$ws.Range("AT1").Value = "Week 45"

# New "day-after" inputs for Fall 22 Week 11
$ws.Range("AT2").Value = 4
$ws.Range("AT6").Value = 5.5
$ws.Range("AT7").Value = 5.5
$ws.Range("AT9").Value = 1.5
$ws.Range("AT10").Value = 2.5

# Update the active selection to reflect where the user ended up after entry
$ws.Activate()
$ws.Range("AS10").Select()
